$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.071.83"
$ws.Range("E2").Value = "  -0.05%  "
$ws.Range("D3").Value = "1.660.46"
$ws.Range("E3").Value = "  -0.35%  "
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").Value = "`'207.89"
$ws.Range("E5").Value = "  -0.80%  "
$ws.Range("D6").Value = "`'0.5170"
$ws.Range("E6").Value = "  -0.22%  "
$ws.Range("D7").Value = "`'1.004"
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("E8").Value = "  -2.13%  "
$ws.Range("D9").Value = "`'0.06300"
$ws.Range("D10").Value = "`'20.89"
$ws.Range("E10").Value = "  -0.18%  "
$ws.Range("D11").Value = "`'0.07543"
$ws.Range("E11").Value = "  +0.70%  "
$ws.Range("D12").Value = "1.660.11"
$ws.Range("E12").Value = "  -1.35%  "
$ws.Range("D13").Value = "`'4.398"
$ws.Range("E13").Value = "  -0.70%  "
$ws.Range("D14").Value = "`'0.5376"
$ws.Range("E14").Value = "  -3.76%  "
$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D15").Value = "0.0₅7958"
$ws.Range("E15").Value = "  +0.75%  "
$ws.Range("B16").Value = "Litecoin"
$ws.Range("C16").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D16").Value = "`'66.14"
$ws.Range("E16").Value = "  +1.01%  "
$ws.Range("D17").Value = "26.080.43"
$ws.Range("E17").Value = "  -0.20%  "
$ws.Range("E18").Value = "  -0.18%  "
$ws.Range("D19").Value = "`'4.696"
$ws.Range("E19").Value = "  -1.98%  "
$ws.Range("D20").Value = "`'187.25"
$ws.Range("E20").Value = "  +0.86%  "
$ws.Range("D21").Value = "`'10.15"
$ws.Range("E21").Value = "  -2.28%  "
$ws.Range("D22").Value = "`'6.190"
$ws.Range("E22").Value = "  +0.55%  "
$ws.Range("E23").Value = "  -0.15%  "
$ws.Range("D24").Value = "`'148.41"
$ws.Range("E24").Value = "  +1.54%  "
$ws.Range("D25").Value = "`'0.1213"
$ws.Range("E25").Value = "  -2.78%  "
$ws.Range("D26").Value = "`'7.382"
$ws.Range("E26").Value = "  -2.18%  "
$ws.Range("D27").Value = "`'15.64"
$ws.Range("E27").Value = "  -0.51%  "
$ws.Range("D28").Value = "`'1.392"
$ws.Range("E28").Value = "  +3.57%  "
$ws.Range("D29").Value = "`'0.05963"
$ws.Range("E29").Value = "  -6.14%  "
$ws.Range("D30").Value = "`'1.261"
$ws.Range("E30").Value = "  -0.92%  "
$ws.Range("D31").Value = "`'3.470"
$ws.Range("E31").Value = "  -0.31%  "
$ws.Range("D32").Value = "`'3.395"
$ws.Range("E32").Value = "  -1.07%  "
$ws.Range("D33").Value = "`'1.636"
$ws.Range("E33").Value = "  +1.01%  "
$ws.Range("D34").Value = "`'0.9841"
$ws.Range("E34").Value = "  -0.97%  "
$ws.Range("D35").Value = "`'2.762"
$ws.Range("E35").Value = "  +2.14%  "
$ws.Range("D36").Value = "`'2.391"
$ws.Range("E36").Value = "  -0.76%  "
$ws.Range("D37").Value = "`'0.5878"
$ws.Range("E37").Value = "  -2.41%  "
$ws.Range("D38").Value = "1.102.76"
$ws.Range("E38").Value = "  +0.93%  "
$ws.Range("D39").Value = "`'0.01597"
$ws.Range("E39").Value = "  -0.61%  "
$ws.Range("D40").Value = "`'5.962"
$ws.Range("E40").Value = "  -1.82%  "
$ws.Range("D41").Value = "`'0.8488"
$ws.Range("E41").Value = "  -1.36%  "
$ws.Range("E42").Value = "  -0.09%  "
$ws.Range("D43").Value = "`'99.86"
$ws.Range("E43").Value = "  +0.39%  "
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("E45").Value = "  +2.36%  "
$ws.Range("D46").Value = "`'55.03"
$ws.Range("E46").Value = "  -1.89%  "
$ws.Range("E47").Value = "  -0.65%  "
$ws.Range("E48").Value = "  +1.03%  "
$ws.Range("D49").Value = "`'0.05231"
$ws.Range("E49").Value = "  -0.30%  "
$ws.Range("D50").Value = "`'0.4242"
$ws.Range("E50").Value = "  -0.48%  "
$ws.Range("E51").Value = "  -0.69%  "
